$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $val) {
    # Force text interpretation so numeric-looking strings (e.g. "68.590.67",
    # "1.00", "  +0.68%  ") are not auto-converted to numbers/dates by Excel,
    # then restore the cell's style so no stray number-format style sticks
    # around (matches the original unstyled inline-string cells).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-Cell "D2" "68.590.67"
Set-Cell "E2" "  +0.68%  "

# Row 3 - Ethereum
Set-Cell "D3" "3.762.27"
Set-Cell "E3" "  -0.73%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  +0.11%  "

# Row 5 - BNB
Set-Cell "D5" "593.74"
Set-Cell "E5" "  -0.65%  "

# Row 6 - Solana
Set-Cell "D6" "167.26"

# Row 7 - LidoStakedEther
Set-Cell "D7" "3.759.74"
Set-Cell "E7" "  -0.78%  "

# Row 9 - XRP
Set-Cell "E9" "  -1.06%  "

# Row 10 - Dogecoin
Set-Cell "E10" "  -2.98%  "

# Row 11 - Toncoin
Set-Cell "D11" "6.41"
Set-Cell "E11" "  -1.81%  "

# Row 12 - Cardano
Set-Cell "D12" "0.450"
Set-Cell "E12" "  -1.03%  "

# Row 13 - ShibaInu
Set-Cell "E13" "  -7.27%  "

# Row 14 - Avalanche
Set-Cell "E14" "  -1.86%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Cell "D15" "4.395.16"
Set-Cell "E15" "  -0.61%  "

# Row 16 - WrappedEther
Set-Cell "D16" "3.767.76"
Set-Cell "E16" "  -0.75%  "

# Row 17 - WrappedBTC
Set-Cell "D17" "68.560.33"
Set-Cell "E17" "  +0.86%  "

# Row 18 - Chainlink
Set-Cell "E18" "  -4.24%  "

# Row 19 - TRON
Set-Cell "E19" "  +0.64%  "

# Row 21 - Uniswap
Set-Cell "D21" "10.75"
Set-Cell "E21" "  +1.37%  "

# Row 22 - BitcoinCash
Set-Cell "D22" "465.70"
Set-Cell "E22" "  -0.83%  "

# Row 23 - Polygon
Set-Cell "D23" "0.698"
Set-Cell "E23" "  -3.32%  "

# Row 24 & 25 swap - PEPE <-> Litecoin
Set-Cell "B24" "Litecoin"
Set-Cell "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-Cell "D24" "84.22"
Set-Cell "E24" "  +0.53%  "

Set-Cell "B25" "PEPE"
Set-Cell "C25" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-Cell "D25" "0.0000147"
Set-Cell "E25" "  -2.28%  "

# Row 26 - Fetch.AI
Set-Cell "E26" "  -2.96%  "

# Row 27 - InternetComputer(DFINITY)
Set-Cell "D27" "11.96"
Set-Cell "E27" "  -1.74%  "

# Row 28 - RenderToken
Set-Cell "D28" "10.05"
Set-Cell "E28" "  -3.92%  "

# Row 29 - Dai
Set-Cell "E29" "  -0.10%  "

# Row 30 - WrappedeETH
Set-Cell "D30" "3.910.08"
Set-Cell "E30" "  -0.66%  "

# Row 31 - PancakeSwap
Set-Cell "E31" "  -4.90%  "

# Row 32 - NEARProtocol
Set-Cell "E32" "  -3.75%  "

# Row 33 - EthereumClassic
Set-Cell "D33" "30.10"
Set-Cell "E33" "  -1.74%  "

# Row 34 - ImmutableX
Set-Cell "E34" "  -3.19%  "

# Row 35 - Aptos
Set-Cell "D35" "9.22"

# Row 36 - Binance-PegBSC-USD
Set-Cell "D36" "0.999"

# Row 37 - RenzoRestakedETH
Set-Cell "D37" "3.715.15"
Set-Cell "E37" "  -0.85%  "

# Row 38 - Hedera
Set-Cell "E38" "  -3.79%  "

# Row 39 - dogwifhat
Set-Cell "E39" "  -9.19%  "

# Row 40 - Kaspa
Set-Cell "E40" "  -1.12%  "

# Row 41 - Mantle
Set-Cell "D41" "1.00"
Set-Cell "E41" "  -0.78%  "

# Row 42 - Filecoin
Set-Cell "E42" "  -1.17%  "

# Row 43 - FirstDigitalUSD
Set-Cell "E43" "  +0.05%  "

# Row 44 - USDe
Set-Cell "E44" "  -0.02%  "

# Row 45 & 46 swap - TheGraph <-> Arweave
Set-Cell "B45" "Arweave"
Set-Cell "C45" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-Cell "D45" "44.04"
Set-Cell "E45" "  +8.53%  "

Set-Cell "B46" "TheGraph"
Set-Cell "C46" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-Cell "D46" "0.303"
Set-Cell "E46" "  -3.79%  "

# Row 47 - OKB
Set-Cell "D47" "46.97"
Set-Cell "E47" "  +2.93%  "

# Row 48 - Stacks
Set-Cell "E48" "  -2.10%  "

# Row 49 - Cosmos
Set-Cell "D49" "8.50"
Set-Cell "E49" "  -2.45%  "

# Row 50 - Monero
Set-Cell "D50" "145.13"
Set-Cell "E50" "  +2.09%  "

# Row 51 - Bittensor
Set-Cell "D51" "389.86"
Set-Cell "E51" "  -3.07%  "
